# "Add files via upload" — extend the fielding-stats table on 工作表1
# with three more historical seasons (2016, 2015, 2014) in column A,
# then leave the selection where the author left it (B9).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(7, 1).Value = 2016
$ws.Cells.Item(8, 1).Value = 2015
$ws.Cells.Item(9, 1).Value = 2014

$ws.Range("B9").Select()
